# Care Connect.pptx edit script
#
# 1. The cached "datetimeFigureOut" field text on the slide master and on
#    every one of the 11 slide layouts advances one day: 1/7/2024 -> 1/8/2024.
# 2. Slide 1 ("Presented by : Samita Maggo") gets an extra space inserted
#    after the colon: ": " -> ":  ".
# 3. Slide 4 (Tools & Technologies / Thymeleaf line) gets an extra space
#    inserted in the leading whitespace run: "     " (5) -> "      " (6).

$p = $ppt.ActivePresentation

# --- 1. Date placeholder on the slide master ---------------------------
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "1/8/2024"
    }
}

# --- ... and on every slide layout --------------------------------------
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "1/8/2024"
        }
    }
}

# --- 2. Slide 1: "Presented by : Samita Maggo" -> extra space after ":" -
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange
$full1 = $tr1.Text
$colonIdx = $full1.IndexOf(": ")
if ($colonIdx -ge 0) {
    $ch1 = $tr1.Characters($colonIdx + 1, 2)
    $ch1.Text = ":  "
}

# --- 3. Slide 4: "     Thymeleaf" -> one more leading space -------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange
$full4 = $tr4.Text
$spacesIdx = $full4.IndexOf("     Thymeleaf")
if ($spacesIdx -ge 0) {
    $ch4 = $tr4.Characters($spacesIdx + 1, 5)
    $ch4.Text = "      "
}

Write-Output "edit.ps1 applied"
